$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 18:59"

# --- Country name changes caused by the re-sort (column A) ---
$countryNames = @{
    20 = "Turquia"
    21 = "Italia"
    84 = "Libano"
    85 = "Paraguay"
    170 = "Birmania"
    171 = "Taiwan"
    214 = "Montserrat"
    215 = "Islas Malvinas"
}
foreach ($row in $countryNames.Keys) {
    $ws.Cells.Item($row, 1).Value = $countryNames[$row]
}

# --- Updated statistics (columns B..H) for the affected rows ---
$statRows = @{
    4 = @(5924741, 9111, 3220168, 2523088, 0, 371, 181485)
    5 = @(3636167, 8950, 2778709, 741812, 0, 195, 115646)
    6 = @(3211848, 46967, 2445975, 706568, 0, 759, 59305)
    12 = @(423224, 2415, 0, 0, 0, 52, 28924)
    16 = @(327798, 1184, 0, 0, 0, 16, 41449)
    20 = @(261194, 1502, 238795, 16236, 0, 24, 6163)
    21 = @(261174, 876, 206015, 19714, 0, 4, 35445)
    23 = @(236810, 693, 209600, 17872, 0, 2, 9338)
    27 = @(125810, 163, 111862, 4862, 0, 3, 9086)
    28 = @(117498, 232, 114318, 2986, 0, 0, 194)
    32 = @(106245, 1773, 83810, 21577, 0, 11, 858)
    44 = @(69651, 1118, 57891, 9130, 0, 19, 2630)
    54 = @(49719, 0, 46311, 3222, 0, 1, 186)
    74 = @(22414, 233, 16650, 5348, 0, 1, 416)
    84 = @(13687, 532, 3723, 9826, 0, 12, 138)
    85 = @(13602, 0, 7649, 5734, 0, 0, 219)
    95 = @(8904, 29, 8399, 449, 0, 0, 56)
    108 = @(5383, 0, 5282, 41, 0, 0, 60)
    157 = @(1184, 85, 178, 991, 0, 0, 15)
    168 = @(710, 6, 658, 10, 0, 0, 42)
    170 = @(504, 30, 341, 157, 0, 0, 6)
    171 = @(487, 0, 457, 23, 0, 0, 7)
    214 = @(13, 0, 12, 0, 0, 0, 1)
    215 = @(13, 0, 13, 0, 0, 0, 0)
}
foreach ($row in $statRows.Keys) {
    $vals = $statRows[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $vals[$i]
    }
}
